$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 82.07692
$ws.Range("I9").Value = 60.636364
$ws.Range("J9").Value = 200
$ws.Range("K9").Value = 60.636364
$ws.Range("L9").Value = 200
$ws.Range("M9").Value = 108.363636
$ws.Range("N9").Value = -538

$ws.Range("H112").Value = 1714.2858
$ws.Range("I112").Value = 1000.3333
$ws.Range("J112").Value = 2249.75
$ws.Range("K112").Value = 3000.9999
$ws.Range("L112").Value = 6749.25
$ws.Range("M112").Value = -1892.9999
$ws.Range("N112").Value = -8965.25

$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("M116").ClearContents()
$ws.Range("N116").ClearContents()

$ws.Range("H127").Value = 1747.75
$ws.Range("J127").Value = 4217
$ws.Range("L127").Value = 12651
$ws.Range("N127").Value = -22571

$ws.Range("H129").Value = 19484.686
$ws.Range("I129").Value = 691.3333
$ws.Range("J129").Value = 20590.176
$ws.Range("K129").Value = 2073.9999
$ws.Range("L129").Value = 61770.528
$ws.Range("M129").Value = 2926.0001
$ws.Range("N129").Value = -71770.52799999999

$ws.Range("H137").Value = 1376.3636
$ws.Range("I137").Value = 1167.5
$ws.Range("J137").Value = 1933.3334
$ws.Range("K137").Value = 3502.5
$ws.Range("L137").Value = 5800.0002
$ws.Range("M137").Value = -952.5
$ws.Range("N137").Value = -10900.0002

$ws.Range("H138").Value = 1904.6136
$ws.Range("J138").Value = 2059.6775
$ws.Range("L138").Value = 6179.032499999999
$ws.Range("N138").Value = -16459.0325

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19121.295
$ws.Range("I32").Value = 20697.11
$ws.Range("J32").Value = 4676.3335
$ws.Range("K32").Value = 20697.11
$ws.Range("L32").Value = 4676.3335
$ws.Range("M32").Value = -20410.11
$ws.Range("N32").Value = -5250.3335

$ws.Range("H61").Value = 2249.9473
$ws.Range("I61").Value = 1767.7858
$ws.Range("J61").Value = 3600
$ws.Range("K61").Value = 1767.7858
$ws.Range("L61").Value = 3600
$ws.Range("M61").Value = -1555.7858
$ws.Range("N61").Value = -4024

$ws.Range("H64").Value = 30055.5
$ws.Range("I64").Value = 25555.5
$ws.Range("J64").Value = 34555.5
$ws.Range("K64").Value = 25555.5
$ws.Range("L64").Value = 34555.5
$ws.Range("M64").Value = -25307.5
$ws.Range("N64").Value = -35051.5

$ws.Range("H67").Value = 30055.5
$ws.Range("I67").Value = 25555.5
$ws.Range("J67").Value = 34555.5
$ws.Range("K67").Value = 25555.5
$ws.Range("L67").Value = 34555.5
$ws.Range("M67").Value = -24697.5
$ws.Range("N67").Value = -36271.5

$ws.Range("H74").Value = 3278
$ws.Range("I74").Value = 3472.5
$ws.Range("J74").Value = 2500
$ws.Range("K74").Value = 3472.5
$ws.Range("L74").Value = 2500
$ws.Range("M74").Value = -2598.5
$ws.Range("N74").Value = -4248

$ws.Range("H77").Value = 3278
$ws.Range("I77").Value = 3472.5
$ws.Range("J77").Value = 2500
$ws.Range("K77").Value = 17362.5
$ws.Range("L77").Value = 12500
$ws.Range("M77").Value = -12994.5
$ws.Range("N77").Value = -21236

$ws.Range("H102").Value = 1949.1666
$ws.Range("I102").Value = 1939
$ws.Range("J102").Value = 2000
$ws.Range("K102").Value = 1939
$ws.Range("L102").Value = 2000
$ws.Range("M102").Value = -317
$ws.Range("N102").Value = -5244

$ws.Range("H122").Value = 3701.625
$ws.Range("I122").Value = 3323.6
$ws.Range("K122").Value = 9970.799999999999
$ws.Range("M122").Value = -7520.799999999999

$ws.Range("H132").Value = 1855.4584
$ws.Range("I132").Value = 1319.0625
$ws.Range("J132").Value = 2928.25
$ws.Range("K132").Value = 3957.1875
$ws.Range("L132").Value = 8784.75
$ws.Range("M132").Value = -1427.1875
$ws.Range("N132").Value = -13844.75

$ws.Range("H136").Value = 2249.9473
$ws.Range("I136").Value = 1767.7858
$ws.Range("J136").Value = 3600
$ws.Range("K136").Value = 5303.357400000001
$ws.Range("L136").Value = 10800
$ws.Range("M136").Value = -2753.357400000001
$ws.Range("N136").Value = -15900

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H103").Value = 38333.332
$ws.Range("J103").Value = 38333.332
$ws.Range("L103").Value = 38333.332
$ws.Range("N103").Value = -40677.332

$ws.Range("H105").Value = 2182.5454
$ws.Range("I105").Value = 2113
$ws.Range("J105").Value = 2495.5
$ws.Range("K105").Value = 2113
$ws.Range("L105").Value = 2495.5
$ws.Range("M105").Value = -366
$ws.Range("N105").Value = -5989.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 12500
$ws.Range("J4").Value = 20000
$ws.Range("L4").Value = 20000
$ws.Range("N4").Value = -20224

$ws.Range("H7").Value = 61.714287
$ws.Range("I7").Value = 48.2
$ws.Range("J7").Value = 69.22221999999999
$ws.Range("K7").Value = 48.2
$ws.Range("L7").Value = 69.22221999999999
$ws.Range("M7").Value = 64.8
$ws.Range("N7").Value = -295.22222

$ws.Range("H31").Value = 3638280.5
$ws.Range("I31").Value = 1938.8334
$ws.Range("J31").Value = 15386462
$ws.Range("K31").Value = 1938.8334
$ws.Range("L31").Value = 15386462
$ws.Range("M31").Value = -1643.8334
$ws.Range("N31").Value = -15387052

$ws.Range("H34").Value = 3638280.5
$ws.Range("I34").Value = 1938.8334
$ws.Range("J34").Value = 15386462
$ws.Range("K34").Value = 1938.8334
$ws.Range("L34").Value = 15386462
$ws.Range("M34").Value = -1736.8334
$ws.Range("N34").Value = -15386866

$ws.Range("H58").Value = 1172.7273
$ws.Range("I58").Value = 842.8570999999999
$ws.Range("K58").Value = 842.8570999999999
$ws.Range("M58").Value = -639.8570999999999

$ws.Range("H105").Value = 764.5454999999999
$ws.Range("I105").Value = 506.25
$ws.Range("J105").Value = 1453.3334
$ws.Range("K105").Value = 506.25
$ws.Range("L105").Value = 1453.3334
$ws.Range("M105").Value = 1240.75
$ws.Range("N105").Value = -4947.3334

$ws.Range("H132").Value = 3288.1177
$ws.Range("I132").Value = 2290
$ws.Range("J132").Value = 4714
$ws.Range("K132").Value = 6870
$ws.Range("L132").Value = 14142
$ws.Range("M132").Value = -4340
$ws.Range("N132").Value = -19202

$ws.Range("H134").Value = 985.129
$ws.Range("I134").Value = 969.8276
$ws.Range("J134").Value = 1207
$ws.Range("K134").Value = 2909.4828
$ws.Range("L134").Value = 3621
$ws.Range("M134").Value = -374.4827999999998
$ws.Range("N134").Value = -8691

$ws.Range("H136").Value = 1172.7273
$ws.Range("I136").Value = 842.8570999999999
$ws.Range("K136").Value = 2528.5713
$ws.Range("M136").Value = 21.42870000000039

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 34.04762
$ws.Range("I2").Value = 14.272727
$ws.Range("J2").Value = 55.8
$ws.Range("K2").Value = 14.272727
$ws.Range("L2").Value = 55.8
$ws.Range("M2").Value = 98.727273
$ws.Range("N2").Value = -281.8

$ws.Range("H58").Value = 1438090.8
$ws.Range("J58").Value = 10251.538
$ws.Range("L58").Value = 10251.538
$ws.Range("N58").Value = -10805.538

$ws.Range("H122").Value = 2433.7334
$ws.Range("I122").Value = 1401.2
$ws.Range("J122").Value = 2950
$ws.Range("K122").Value = 4203.6
$ws.Range("L122").Value = 8850
$ws.Range("M122").Value = -1753.6
$ws.Range("N122").Value = -13750

$ws.Range("H126").Value = 8341583
$ws.Range("I126").Value = 16000
$ws.Range("J126").Value = 16667166
$ws.Range("K126").Value = 48000
$ws.Range("L126").Value = 50001498
$ws.Range("M126").Value = -45530
$ws.Range("N126").Value = -50006438

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H133").Value = 26013
$ws.Range("J133").Value = 26013
$ws.Range("L133").Value = 26013
$ws.Range("N133").Value = -31073

$ws.Range("H136").Value = 9114
$ws.Range("I136").Value = 14037.125
$ws.Range("J136").Value = 2549.8333
$ws.Range("K136").Value = 42111.375
$ws.Range("L136").Value = 7649.499899999999
$ws.Range("M136").Value = -39561.375
$ws.Range("N136").Value = -12749.4999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H76").Value = 20293.25
$ws.Range("J76").Value = 24173
$ws.Range("L76").Value = 24173
$ws.Range("N76").Value = -24803

$ws.Range("H79").Value = 20293.25
$ws.Range("J79").Value = 24173
$ws.Range("L79").Value = 24173
$ws.Range("N79").Value = -26357

$ws.Range("H122").Value = 2187.2354
$ws.Range("I122").Value = 2670.8
$ws.Range("J122").Value = 1496.4286
$ws.Range("K122").Value = 8012.400000000001
$ws.Range("L122").Value = 4489.2858
$ws.Range("M122").Value = -5562.400000000001
$ws.Range("N122").Value = -9389.2858

$ws.Range("H136").Value = 14317.333
$ws.Range("I136").Value = 15994.5
$ws.Range("J136").Value = 900
$ws.Range("K136").Value = 47983.5
$ws.Range("L136").Value = 2700
$ws.Range("M136").Value = -45433.5
$ws.Range("N136").Value = -7800
